# Commit: "Updated most important problem list with 10 new responses."
#
# The underlying survey-response table on "Sheet 2" gained coded values
# (column C = numeric category code, column D = category label) for the
# tail rows (108-115) that previously only had the raw free-text answer
# (column A) plus the constant E/F columns. This mirrors a researcher
# having gone back and manually coded the newest open-ended responses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 2")

# row -> (category code, category label)
$ws.Range("C108").Value = 6
$ws.Range("D108").Value = "Mental health"

$ws.Range("C109").Value = 2
$ws.Range("D109").Value = "Public health messaging / gov't handing of COVID"

$ws.Range("C110").Value = 5
$ws.Range("D110").Value = "Chronic disease (cancer, heart disease, etc)"

$ws.Range("C111").Value = 7
$ws.Range("D111").Value = "Access to housing and food"

$ws.Range("C112").Value = 7
$ws.Range("D112").Value = "Access to housing and food"

$ws.Range("C113").Value = 7
$ws.Range("D113").Value = "Access to housing and food"

$ws.Range("C114").Value = 2
$ws.Range("D114").Value = "Mental health"

$ws.Range("C115").Value = 9
$ws.Range("D115").Value = "Inequality and discrimination"

# Match the author's final on-screen selection/scroll position captured
# in the saved file (cell D114 selected, sheet scrolled so B95 is the
# top-left visible cell).
$ws.Range("D114").Select()
